$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.3396209048899, 0.7817619268305633, 0.9554726193415759, 1.099153532933269, 0.998)
    3  = @(6.893910660182505, 0.9015149716923674, 1.960136467179636, 1.099153532933269, 0.998)
    4  = @(3.642976410215606, 0.8200517269391167, 1.538269377134317, 1.099153532933269, 0.998)
    5  = @(2.997774918752558, 0.998030937379355, 1.260670400900648, 1.099153532933269, 0.998)
    6  = @(2.103375112193536, 0.9770801135847219, 1.152925594936625, 1.099153532933269, 0.998)
    7  = @(1.782945229375743, 0.9990003583255496, 1.050046863520349, 1.099153532933269, 0.998)
    8  = @(2.200748730901851, 0.9974791451534021, 1.232858037995043, 1.099153532933269, 0.998)
    9  = @(15.54654381901065, 0.8141518114531443, 3.250987742369093, 1.099153532933269, 0.998)
    10 = @(1.562861492608981, 0.9953411582976511, 1.008048432933262, 1.099153532933269, 0.998)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("G$row").Value = $vals[3]
    $ws.Range("H$row").Value = $vals[4]
}
